# v2.0 Fix push buttons
# The "Enclosure" row's "Where to Buy" source changed from "Computer City"
# to "Micro Center".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C15").Value = "Micro Center"

# Leave the selection where the user last edited, scrolled down a bit so
# row 15 is in view (best-effort; harmless if unsupported).
try {
    $excel.ActiveWindow.ScrollRow = 12
} catch {}

$ws.Range("C15").Select()
